{"js": "// Change the trailing \"1\" in \"Testing word doc 1\" to \"2\". In the author's\n// edit the run that used to read \"Testing word doc 1\" got split into two\n// runs with identical formatting: \"Testing word doc \" and \"2\".\n\nconst body = context.document.body;\n\n// Locate the sentence we need to edit (narrow search so we don't touch an\n// unrelated \"1\" elsewhere in the document).\nconst sentenceHits = body.search(\"Testing word doc 1\", { matchCase: true });\nsentenceHits.load(\"items\");\nawait context.sync();\n\nif (sentenceHits.items.length === 0) {\n  throw new Error('Could not find \"Testing word doc 1\" in the document body.');\n}\nconst sentence = sentenceHits.items[0];\n\n// Read that run's own OOXML so the freshly split-off run keeps the same\n// <w:rPr> (here: xml:lang=\"en-US\") as the text it came from.\nconst ooxml = sentence.getOoxml();\nawait context.sync();\nconst langMatch = /<w:lang\\b[^>]*\\/>/.exec(ooxml.value);\nconst langTag = langMatch ? langMatch[0] : '<w:lang w:val=\"en-US\"/>';\n\n// Within that sentence, find the trailing \"1\" character to replace.\nconst digitHits = sentence.search(\"1\", { matchCase: true });\ndigitHits.load(\"items\");\nawait context.sync();\nif (digitHits.items.length === 0) {\n  throw new Error('Could not find the trailing \"1\" to replace.');\n}\nconst digit = digitHits.items[digitHits.items.length - 1];\n\n// Office.js's insertOoxml requires the flat-OPC envelope; inside it we\n// supply a single new run (\"2\") that replaces the \"1\" range. Because the\n// replacement is a distinct <w:r> (not a text patch inside the existing\n// run), the result is two sibling runs: \"Testing word doc \" + \"2\".\nconst flatOpc =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body><w:p>' +\n  '<w:r><w:rPr>' + langTag + '</w:rPr><w:t>2</w:t></w:r>' +\n  '</w:p></w:body>' +\n  '</w:document>' +\n  '</pkg:xmlData></pkg:part></pkg:package>';\n\ndigit.insertOoxml(flatOpc, \"Replace\");\nawait context.sync();\n", "ps1": "# Change the trailing \"1\" in \"Testing word doc 1\" to \"2\". In the author's\n# edit the run that used to read \"Testing word doc 1\" got split into two\n# runs with identical formatting: \"Testing word doc \" and \"2\".\n\n$d = $word.ActiveDocument\n\n# Locate the sentence we need to edit (narrow search so we don't touch an\n# unrelated \"1\" elsewhere in the document).\n$sentence = $d.Content\n$found = $sentence.Find.Execute(\"Testing word doc 1\")\nif (-not $found) {\n    throw 'Could not find \"Testing word doc 1\" in the document body.'\n}\n\n# Within that sentence, narrow down to just the trailing \"1\" character.\n$digit = $d.Range($sentence.End - 1, $sentence.End)\nif ($digit.Text -ne \"1\") {\n    throw 'Could not find the trailing \"1\" to replace.'\n}\n\n# Replace just that \"1\" range with a brand-new <w:r> (\"2\") via InsertXML,\n# wrapped in the flat-OPC envelope. Because the replacement is a distinct\n# run (not a text patch inside the existing run), the result is two\n# sibling runs: \"Testing word doc \" + \"2\", both tagged xml:lang=\"en-US\"\n# like the original run.\n$flatOpc = @\"\n<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t>2</w:t></w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>\n\"@\n\n$digit.InsertXML($flatOpc)\n"}
